$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the external workbook link (homologación de precios de insumos) ---
$links = $wb.LinkSources()
if ($links) {
    foreach ($l in $links) {
        $wb.BreakLink($l, 1)
    }
}

# --- Header row: drop the tall/wrapped row height, selection now on B1 ---
$ws.Rows.Item(1).AutoFit()
$ws.Range("B1").Select()

# --- Append new monthly price rows 109-117 ---
$newRows = @(
    @(109, 45261, 168.20156001145301, 159.035415834426,    157.98174249704101, 182.72643848899901),
    @(110, 45292, 168.20156001145301, 159.035415834426,    157.98174249704101, 171.602023972756),
    @(111, 45323, 168.20156001145301, 159.035415834426,    157.98174249704101, 174.91479838073499),
    @(112, 45352, 168.20156001145301, 159.035415834426,    157.98174249704101, 181.495851693942),
    @(113, 45383, 173.00731886892299, 163.67394879626301,  162.58954331987201, 190.96306729767201),
    @(114, 45413, 173.00731886892299, 163.673948796264,    162.58954331987201, 196.322380513197),
    @(115, 45444, 173.00731886892299, 165.137198067111,    164.66692458275401, 194.93070915982),
    @(116, 45474, 173.00731886892299, 181.52968459686701,  180.19666404162999, 189.946264977749),
    @(117, 45505, 173.00731886892299, 183.02580491190699,  181.74241763149001, 188.746081090817)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

Write-Host "done"
